# Update header cells in the mock user data worksheet so that the
# username/password column headers become Robot Framework style
# variable placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = '${username}'
$ws.Range("B1").Value = '${password}'

# Reset the selection back to the default top-left cell so the saved
# worksheet does not retain a stale selection/activeCell range.
$ws.Range("A1").Select()
